# Generation of instances similar to Templemeier 1996
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Generic: rename the NrTimeBucketWithoutUncertainty row into a "Before"
# row and add a new "After" row underneath it; NrResources drops to 3.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Generic")
$ws.Range("B5").Value = 3
$ws.Range("A8").Value = "NrTimeBucketWithoutUncertaintyBefore"
$ws.Range("B8").Value = 0
# Clone A8's label style onto the new A9 row, then fill in the real content.
$ws.Range("A8").Copy($ws.Range("A9"))
$ws.Range("A9").Value = "NrTimeBucketWithoutUncertaintyAfter"
$ws.Range("B9").Value = 3

# ---------------------------------------------------------------------
# Productdata: updated starting inventories / setup costs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("E2").Value = 9.500400000000001
$ws.Range("E3").Value = 3.4752
$ws.Range("C4").Value = 895
$ws.Range("E4").Value = 4.0816
$ws.Range("C5").Value = 895
$ws.Range("E5").Value = 1.700666666666667
$ws.Range("C6").Value = 895
$ws.Range("E6").Value = 3.0612
$ws.Range("C7").Value = 522
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("C9").Value = 150
$ws.Range("E9").Value = 0

# ---------------------------------------------------------------------
# ForecastedAverageDemand: updated demand forecasts
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("G2").Value = 272
$ws.Range("H2").Value = 43
$ws.Range("I2").Value = 74
$ws.Range("G3").Value = 246
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 75
$ws.Range("G4").Value = 287
$ws.Range("H4").Value = 47
$ws.Range("I4").Value = 74
$ws.Range("G5").Value = 259
$ws.Range("H5").Value = 44
$ws.Range("I5").Value = 76
$ws.Range("G6").Value = 268
$ws.Range("H6").Value = 46
$ws.Range("I6").Value = 78
$ws.Range("G7").Value = 227
$ws.Range("H7").Value = 43
$ws.Range("I7").Value = 79

# ---------------------------------------------------------------------
# ForcastedStandardDeviation: updated standard deviation forecasts
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("G2").Value = 6.799999999999999
$ws.Range("H2").Value = 1.075
$ws.Range("I2").Value = 1.85
$ws.Range("G3").Value = 11.685
$ws.Range("H3").Value = 2.137499999999999
$ws.Range("I3").Value = 3.562499999999999
$ws.Range("G4").Value = 19.44424999999999
$ws.Range("H4").Value = 3.184249999999999
$ws.Range("I4").Value = 5.013499999999999
$ws.Range("G5").Value = 22.267525
$ws.Range("H5").Value = 3.7829
$ws.Range("I5").Value = 6.5341
$ws.Range("G6").Value = 27.43716999999999
$ws.Range("H6").Value = 4.709364999999999
$ws.Range("I6").Value = 7.985444999999999
$ws.Range("G7").Value = 26.59072325
$ws.Range("H7").Value = 5.03700925
$ws.Range("I7").Value = 9.254040249999999

# ---------------------------------------------------------------------
# Capacity: only 3 time buckets remain (rows 5:9 removed), values updated
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("A5:B9").EntireRow.Delete()
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 8503.333333333332
$ws.Range("B4").Value = 45067.66666666666

# ---------------------------------------------------------------------
# ProcessingTime: only 3 time buckets remain (columns E:I removed),
# values updated for the remaining buckets
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("E1:I9").EntireColumn.Delete()
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 10
$ws.Range("D4").Value = 28
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 10
